$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.142.75'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '''1.655.70'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").Value = '''218.28'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '''0.5286'
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("D8").Value = '''0.2608'
$ws.Range("E8").Value = '  -2.03%  '
$ws.Range("D9").Value = '''0.06345'
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("D10").Value = '''20.47'
$ws.Range("E10").Value = '  -2.36%  '
$ws.Range("D11").Value = '''0.07781'
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").Value = '''4.498'
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").Value = '''1.650.05'
$ws.Range("E13").Value = '  -1.47%  '
$ws.Range("D14").Value = '''0.5479'
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("D15").Value = '''0.0₅8171'
$ws.Range("E15").Value = '  -0.92%  '
$ws.Range("D16").Value = '''65.42'
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("D17").Value = '''26.131.28'
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = '''4.563'
$ws.Range("E19").Value = '  -2.11%  '
$ws.Range("D20").Value = '''192.79'
$ws.Range("E20").Value = '  -1.05%  '
$ws.Range("D21").Value = '''10.08'
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").Value = '''6.033'
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("D23").Value = '''1.003'
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").Value = '''141.95'
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("D25").Value = '''0.1250'
$ws.Range("E25").Value = '  +0.86%  '
$ws.Range("D26").Value = '''7.274'
$ws.Range("E26").Value = '  +1.10%  '
$ws.Range("D27").Value = '''16.19'
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").Value = '''1.435'
$ws.Range("E28").Value = '  +1.42%  '
$ws.Range("D29").Value = '''0.05939'
$ws.Range("E29").Value = '  -3.89%  '
$ws.Range("D30").Value = '''1.284'
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("D31").Value = '''3.522'
$ws.Range("E31").Value = '  -2.10%  '
$ws.Range("D32").Value = '''3.249'
$ws.Range("E32").Value = '  -1.42%  '
$ws.Range("D33").Value = '''1.573'
$ws.Range("E33").Value = '  -3.63%  '
$ws.Range("E34").Value = '  +0.79%  '
$ws.Range("D35").Value = '''0.9504'
$ws.Range("E35").Value = '  -2.39%  '
$ws.Range("D36").Value = '''2.408'
$ws.Range("E36").Value = '  -0.78%  '
$ws.Range("D37").Value = '''0.5655'
$ws.Range("E37").Value = '  -1.92%  '
$ws.Range("D38").Value = '''0.01610'
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("D39").Value = '''5.811'
$ws.Range("D40").Value = '''0.8487'
$ws.Range("E40").Value = '  -1.10%  '
$ws.Range("E41").Value = '  -0.32%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '''102.47'
$ws.Range("E42").Value = '  +2.15%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '''1.024.44'
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("D44").Value = '''1.799.25'
$ws.Range("D45").Value = '''57.21'
$ws.Range("E45").Value = '  -1.00%  '
$ws.Range("D46").Value = '''1.007'
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("D47").Value = '''0.4285'
$ws.Range("E47").Value = '  +1.56%  '
$ws.Range("D48").Value = '''1.473'
$ws.Range("E48").Value = '  -0.47%  '
$ws.Range("D50").Value = '''7.793'
$ws.Range("E50").Value = '  -3.63%  '
$ws.Range("D51").Value = '''0.09720'
$ws.Range("E51").Value = '  -0.71%  '
